$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    @{ Row = 10; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 14; DAMSLTag = 'ba'; DialogAct = 'Appreciation' }
    @{ Row = 18; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 29; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 72; DAMSLTag = 'ba'; DialogAct = 'Appreciation' }
    @{ Row = 76; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 77; DAMSLTag = 'ba'; DialogAct = 'Appreciation' }
    @{ Row = 79; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 80; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 91; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 92; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 99; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 168; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 189; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 192; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 193; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 221; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 222; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 229; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 236; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 249; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 252; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 255; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 278; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 290; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 291; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 315; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 320; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 343; DAMSLTag = 'ba'; DialogAct = 'Appreciation' }
    @{ Row = 349; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 378; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 382; DAMSLTag = '%'; DialogAct = 'Uninterpretable' }
    @{ Row = 388; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 396; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 415; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 417; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 425; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 426; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 430; DAMSLTag = 'ba'; DialogAct = 'Appreciation' }
    @{ Row = 440; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 469; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 470; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 471; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 486; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 494; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 506; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 514; DAMSLTag = '%'; DialogAct = 'Uninterpretable' }
    @{ Row = 531; DAMSLTag = 'ba'; DialogAct = 'Appreciation' }
    @{ Row = 532; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 543; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 553; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 555; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 572; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 573; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 579; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.DAMSLTag
    $ws.Cells.Item($u.Row, 10).Value = $u.DialogAct
}

$wb.Save()
